$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BB values (column 54), one per row 1-83, mirroring/extending column BA
$bbValues = @(
    @{Cell="BB1"; Value=45986},
    @{Cell="BB2"; Value=0.7596024993684409},
    @{Cell="BB3"; Value=5.398981879140436},
    @{Cell="BB4"; Value=-5.469777829091811},
    @{Cell="BB5"; Value=-0.5016308114541062},
    @{Cell="BB6"; Value=-0.6214054032250829},
    @{Cell="BB7"; Value=-1.41027180704431},
    @{Cell="BB8"; Value=-0.7919227738455277},
    @{Cell="BB9"; Value=0.8281781605252121},
    @{Cell="BB10"; Value=0.6333580932638085},
    @{Cell="BB11"; Value=-1.871721946158218},
    @{Cell="BB12"; Value=6.947074551031477},
    @{Cell="BB13"; Value=-0.5538313209713124},
    @{Cell="BB14"; Value=-2.926189405619382},
    @{Cell="BB15"; Value=9.498197854786355},
    @{Cell="BB16"; Value=-0.4994768307962829},
    @{Cell="BB17"; Value=-0.6902935135373696},
    @{Cell="BB18"; Value=1.660952585620919},
    @{Cell="BB19"; Value=0.5318249089792459},
    @{Cell="BB20"; Value=0.6320478207229741},
    @{Cell="BB21"; Value=-0.1274490368921875},
    @{Cell="BB22"; Value=-0.6653180828667331},
    @{Cell="BB23"; Value=-2.768682804727675},
    @{Cell="BB24"; Value=2.972562358782},
    @{Cell="BB25"; Value=1.814509864363558},
    @{Cell="BB26"; Value=0.7290697448430592},
    @{Cell="BB27"; Value=4.434377037602189},
    @{Cell="BB28"; Value=-3.830674548307485},
    @{Cell="BB29"; Value=-0.3546604117736365},
    @{Cell="BB30"; Value=0.5841110884116603},
    @{Cell="BB31"; Value=0.3311457281599957},
    @{Cell="BB32"; Value=-0.983338854947732},
    @{Cell="BB33"; Value=0.0367937104471423},
    @{Cell="BB34"; Value=2.049812641278834},
    @{Cell="BB35"; Value=2.087453263562125},
    @{Cell="BB36"; Value=-1.438164479066856},
    @{Cell="BB37"; Value=-0.02684958208732269},
    @{Cell="BB38"; Value=0.8237845450787233},
    @{Cell="BB39"; Value=2.86126545641612},
    @{Cell="BB40"; Value=0.5069532644027532},
    @{Cell="BB41"; Value=-0.2174083661145261},
    @{Cell="BB42"; Value=0.1045745676356802},
    @{Cell="BB43"; Value=0.4936918743079417},
    @{Cell="BB44"; Value=1.322625689988016},
    @{Cell="BB45"; Value=0.6899772607160202},
    @{Cell="BB46"; Value=1.240850467000882},
    @{Cell="BB47"; Value=2.6},
    @{Cell="BB48"; Value=-0.9},
    @{Cell="BB49"; Value=0.4},
    @{Cell="BB50"; Value=0.1},
    @{Cell="BB51"; Value=4.305581231918552},
    @{Cell="BB52"; Value=-4.290231943195352},
    @{Cell="BB53"; Value=-1.294763502326944},
    @{Cell="BB54"; Value=3.727166552773809},
    @{Cell="BB55"; Value=-0.07597546477697392},
    @{Cell="BB56"; Value=1.620116340063873},
    @{Cell="BB57"; Value=-3.620177172843626},
    @{Cell="BB58"; Value=-0.9596300199833507},
    @{Cell="BB59"; Value=3.509382378526155},
    @{Cell="BB60"; Value=-3.103464835402306},
    @{Cell="BB61"; Value=-0.9738690852321383},
    @{Cell="BB62"; Value=-3.206836580208446},
    @{Cell="BB63"; Value=1.5011358099934},
    @{Cell="BB64"; Value=-0.5164324031961769},
    @{Cell="BB65"; Value=-0.7267585812813877},
    @{Cell="BB66"; Value=-1.863476264789497},
    @{Cell="BB67"; Value=0.5061198669764195},
    @{Cell="BB68"; Value=-1.844723754442185},
    @{Cell="BB69"; Value=-0.3113146714536583},
    @{Cell="BB70"; Value=0.8443262641745406},
    @{Cell="BB71"; Value=0.2255501838065186},
    @{Cell="BB72"; Value=-1.137044400346582},
    @{Cell="BB73"; Value=-0.4792091214565772},
    @{Cell="BB74"; Value=-0.4792091214565772},
    @{Cell="BB75"; Value=-0.4792091214565772},
    @{Cell="BB76"; Value=-0.4792091214565772},
    @{Cell="BB77"; Value=-0.4792091214565772},
    @{Cell="BB78"; Value=-0.4792091214565772},
    @{Cell="BB79"; Value=-0.4792091214565772},
    @{Cell="BB80"; Value=-0.4792091214565772},
    @{Cell="BB81"; Value=-0.4792091214565772},
    @{Cell="BB82"; Value=-0.4792091214565772},
    @{Cell="BB83"; Value=-0.4792091214565772}
)

foreach ($item in $bbValues) {
    $ws.Range($item.Cell).Value = $item.Value
}

# New row 83 in column A (date)
$ws.Range("A83").Value = 46934

# Copy style (bold/border/date-format) from header row BA1 -> BB1
$ws.Range("BA1").Copy() | Out-Null
$ws.Range("BB1").PasteSpecial(-4122) | Out-Null

# Copy style from column A date cell (A82) -> new A83 date cell
$ws.Range("A82").Copy() | Out-Null
$ws.Range("A83").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
